$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: change date separators from "/" to "-" for rows 3..21.
# Some DD-MM-YYYY strings (day <= 12) are ambiguous and Excel's normal
# value-assignment auto-recognizes them as dates, converting them to serial
# numbers. Force those specific cells to Text format first so the literal
# string is preserved, matching the source data (plain text, not a date).
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($r in $dates.Keys) {
    $cell = $ws.Range("A$r")
    $parts = $dates[$r].Split("-")
    $day = [int]$parts[0]
    if ($day -le 12) {
        # Ambiguous as MM-DD-YYYY too - pin the cell to Text so Excel
        # doesn't silently reinterpret the literal as a date serial.
        $cell.NumberFormat = "@"
    }
    $cell.Value = $dates[$r]
}

# Attendance-count updates (D, E, G, H columns) for rows 3-5
$ws.Range("D3").Value = 2
$ws.Range("G3").Value = 2

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0
